$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Replace `old` with `new` the first time it is found anywhere in the
# document. Leaves everything inside one run (a plain text substitution);
# we re-split the run(s) into the desired boundaries afterwards.
function ReplaceOnce($old, $new) {
    $rr = $word.ActiveDocument.Content
    $ok = $rr.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: ReplaceOnce could not find: $old"
    }
}

# Start a new left-to-right run-splitting pass across paragraph $paraIndex.
# Resets the tracked cursor to the start of that paragraph.
function BeginSplit($paraIndex) {
    $p = $word.ActiveDocument.Paragraphs.Item($paraIndex)
    $global:splitCursor = $p.Range.Start
    $global:splitEnd = $p.Range.End
}

# Find `searchText` starting at the tracked cursor (within the paragraph
# being processed) and force a run boundary immediately before it by
# toggling Bold on/off (a no-op formatting change that nonetheless forces
# the engine to materialise a separate <w:r>). Advances the cursor past
# the match so subsequent calls find the next occurrence.
function SplitFrom($searchText) {
    $dd = $word.ActiveDocument
    $rr = $dd.Range($global:splitCursor, $global:splitEnd)
    $found = $rr.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "WARNING: SplitFrom could not find: $searchText"
        return
    }
    $rr.Bold = 1
    $rr.Bold = 0
    $global:splitCursor = $rr.End
}

# ---------------------------------------------------------------------------
# Change 1: paragraph "Description and results from experiment 1" -- split
# the final run into five, inserting a new sentence about stereoacuity.
# ---------------------------------------------------------------------------

$old1 = "led to a response at the fundamental and its second harmonic (6Hz). Surprisingly, no steady-state responses were found for interocular anticorrelation; signal-to-noise ratios at 3Hz were no different from those of the control condition. "
$new1 = "led to a response at the fundamental and its second harmonic (6Hz). We measured stereoacuity in our observers and found that thresholds correlated strongly with the SNRs of the fundamental frequency, but not with its second harmonic. Surprisingly, no steady-state responses were found for interocular anticorrelation; signal-to-noise ratios at 3Hz were no different from those of the control condition. "
ReplaceOnce $old1 $new1

$para1 = 16
BeginSplit $para1
SplitFrom "led to a response at the fundamental and its second harmonic (6Hz)."
SplitFrom " We measured stereoacuity in our observers and found that thresholds correlated strongly with the "
SplitFrom "SNRs"
SplitFrom " of the fundamental frequency, but not with its second harmonic. "
SplitFrom "Surprisingly, no steady-state responses were found for interocular anticorrelation; signal-to-noise ratios at 3Hz were no different from those of the control condition. "

# ---------------------------------------------------------------------------
# Change 2: "Description of a simple model" -> "Description of a generic
# model", split into three runs.
# ---------------------------------------------------------------------------

ReplaceOnce "Description of a simple model" "Description of a generic model"

$para2 = 18
BeginSplit $para2
SplitFrom "Description of a "
SplitFrom "generic"
SplitFrom " model"

# ---------------------------------------------------------------------------
# Change 3: model-description paragraph -- rewrite everything after "our".
# ---------------------------------------------------------------------------

$old3 = " data using an image-based variant of the two-stage contrast gain control model of binocular summation. To generate SSVEPs dependent on interocular correlation, the noise images were filtered with a bank of log-Gabors that had preferred orientations ranging from 0° to 165°, in increments of 15°, and preferred spatial frequencies of 0.5, 1, 2, 4, 8, and 16 cycles/°. The monocular filter responses underwent an early non-linearity and contrast gain control before binocular summation and binocular difference. The sum and difference responses were fed through a second non-linear and contrast gain control. The resulting output was Fourier transformed to generate model SSVEPs. "
$new3 = " data using an image-based variant of the two-stage contrast gain control model of binocular summation. In this model, interocular correlation-dependent SSVEPs were generated from the responses of a bank of disparity-selective log-Gabor filters. The monocular filter responses were fed through a nonlinearity and contrast gain control, followed by binocular summation and binocular difference. The binocular sum and difference responses were fed through a second nonlinearity and binocular contrast gain control, with the output Fourier-transformed to generate model SSVEPs. "
ReplaceOnce $old3 $new3

$para3 = 19
BeginSplit $para3
SplitFrom " data using an image-based variant of the two-stage contrast gain control model of binocular summation. "
SplitFrom "In this model, interocular correlation-dependent SSVEPs were generated from the responses of a bank of disparity-selective log-Gabor filters. The monocular filter responses were fed through a nonlinearity and contrast gain control, followed by binocular summation and binocular difference. The binocular sum and difference responses were fed through a second nonlinearity and binocular contrast gain control"
SplitFrom ","
SplitFrom " with the output Fourier"
SplitFrom "-"
SplitFrom "transformed to generate model SSVEPs. "

# ---------------------------------------------------------------------------
# Change 4: add a "Conclusion" heading (italic) and closing paragraph at the
# end of the document (the last paragraph was empty; a further empty
# paragraph now becomes the heading, and a brand-new paragraph holds the
# concluding sentence).
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertAfter("Conclusion")

$concHeading = $d.Paragraphs.Item($d.Paragraphs.Count)
$concHeading.Range.Italic = 1
$concHeading.Range.ItalicBi = 1

$concHeading.Range.InsertParagraphAfter()

$concBody = $d.Paragraphs.Item($d.Paragraphs.Count)
$concBody.Range.InsertAfter("Neural responses to interocular correlation were measurable using SSVEP methods, yet responses to anticorrelation (differences) were not. ")

Write-Output "Edit complete."
